# Auto-generated edit script: updates Leve profit-calculation workbook cells
# per the scheduled market-data refresh (H/I/J/K/L/M/N columns across ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("H69").Value = 9603.333000000001
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9603.333000000001
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 28809.999
$ws.Range("N69").Value = -30557.999
$ws.Range("H72").Value = 9603.333000000001
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9603.333000000001
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 86429.997
$ws.Range("N72").Value = -95165.997
$ws.Range("H81").Value = 40000
$ws.Range("J81").Value = 40000
$ws.Range("L81").Value = 40000
$ws.Range("N81").Value = -41996
$ws.Range("H84").Value = 40000
$ws.Range("J84").Value = 40000
$ws.Range("L84").Value = 120000
$ws.Range("N84").Value = -129984
$ws.Range("H100").Value = 1940
$ws.Range("I100").Value = 1921.8
$ws.Range("K100").Value = 1921.8
$ws.Range("M100").Value = -1380.8
$ws.Range("H112").Value = 1080.2667
$ws.Range("I112").Value = 600
$ws.Range("J112").Value = 1114.5714
$ws.Range("K112").Value = 1800
$ws.Range("L112").Value = 3343.7142
$ws.Range("M112").Value = -692
$ws.Range("N112").Value = -5559.7142
$ws.Range("M20").ClearContents()
$ws.Range("M35").ClearContents()
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39651.54
$ws.Range("I2").Value = 1103
$ws.Range("K2").Value = 1103
$ws.Range("M2").Value = -990
$ws.Range("H32").Value = 30028.9
$ws.Range("I32").Value = 4806.0376
$ws.Range("J32").Value = 221002
$ws.Range("K32").Value = 4806.0376
$ws.Range("L32").Value = 221002
$ws.Range("M32").Value = -4519.0376
$ws.Range("N32").Value = -221576
$ws.Range("H57").Value = 0
$ws.Range("I57").Value = 0
$ws.Range("K57").Value = 0
$ws.Range("H110").Value = 50055730
$ws.Range("I110").Value = 55617200
$ws.Range("K110").Value = 55617200
$ws.Range("M110").Value = -55615155
$ws.Range("H114").Value = 25266.666
$ws.Range("J114").Value = 25266.666
$ws.Range("L114").Value = 25266.666
$ws.Range("N114").Value = -33944.666
$ws.Range("H116").Value = 39651.54
$ws.Range("I116").Value = 1103
$ws.Range("K116").Value = 1103
$ws.Range("M116").Value = 1191
$ws.Range("M57").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39651.54
$ws.Range("I3").Value = 1103
$ws.Range("K3").Value = 1103
$ws.Range("M3").Value = -989
$ws.Range("H44").Value = 15333.333
$ws.Range("J44").Value = 15333.333
$ws.Range("L44").Value = 15333.333
$ws.Range("N44").Value = -16327.333
$ws.Range("H99").Value = 2283.75
$ws.Range("J99").Value = 2238
$ws.Range("L99").Value = 2238
$ws.Range("N99").Value = -5234
$ws.Range("H107").Value = 55626076
$ws.Range("I107").Value = 111247784
$ws.Range("J107").Value = 4370.3335
$ws.Range("K107").Value = 111247784
$ws.Range("L107").Value = 4370.3335
$ws.Range("M107").Value = -111245864
$ws.Range("N107").Value = -8210.333500000001
$ws.Range("H137").Value = 50780
$ws.Range("J137").Value = 50780
$ws.Range("L137").Value = 50780
$ws.Range("N137").Value = -60980

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19555.105
$ws.Range("I31").Value = 24528.441
$ws.Range("J31").Value = 4279.857
$ws.Range("K31").Value = 24528.441
$ws.Range("L31").Value = 4279.857
$ws.Range("M31").Value = -24233.441
$ws.Range("N31").Value = -4869.857
$ws.Range("H34").Value = 19555.105
$ws.Range("I34").Value = 24528.441
$ws.Range("J34").Value = 4279.857
$ws.Range("K34").Value = 24528.441
$ws.Range("L34").Value = 4279.857
$ws.Range("M34").Value = -24326.441
$ws.Range("N34").Value = -4683.857
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("H58").Value = 7634.5386
$ws.Range("I58").Value = 922.3200000000001
$ws.Range("K58").Value = 922.3200000000001
$ws.Range("M58").Value = -719.3200000000001
$ws.Range("H118").Value = 10000
$ws.Range("J118").Value = 10000
$ws.Range("L118").Value = 10000
$ws.Range("N118").Value = -13314
$ws.Range("H136").Value = 7634.5386
$ws.Range("I136").Value = 922.3200000000001
$ws.Range("K136").Value = 2766.96
$ws.Range("M136").Value = -216.96
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 999.5
$ws.Range("I63").Value = 1000
$ws.Range("J63").Value = 999
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 2997
$ws.Range("M63").Value = -2251
$ws.Range("N63").Value = -4495
$ws.Range("H66").Value = 999.5
$ws.Range("I66").Value = 1000
$ws.Range("J66").Value = 999
$ws.Range("K66").Value = 9000
$ws.Range("L66").Value = 8991
$ws.Range("M66").Value = -5256
$ws.Range("N66").Value = -16479
$ws.Range("H122").Value = 8337.385
$ws.Range("J122").Value = 20961.4
$ws.Range("L122").Value = 188652.6
$ws.Range("N122").Value = -193552.6
$ws.Range("H139").Value = 2375.24
$ws.Range("I139").Value = 1601.4
$ws.Range("J139").Value = 3536
$ws.Range("K139").Value = 4804.200000000001
$ws.Range("L139").Value = 10608
$ws.Range("M139").Value = 335.7999999999993
$ws.Range("N139").Value = -20888

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 7412419.5
$ws.Range("I16").Value = 11455026
$ws.Range("J16").Value = 975.1667
$ws.Range("K16").Value = 11455026
$ws.Range("L16").Value = 975.1667
$ws.Range("M16").Value = -11454856
$ws.Range("N16").Value = -1315.1667
$ws.Range("H133").Value = 53514.285
$ws.Range("J133").Value = 54100
$ws.Range("L133").Value = 54100
$ws.Range("N133").Value = -59160

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 401691.8
$ws.Range("J81").Value = 252114.75
$ws.Range("L81").Value = 504229.5
$ws.Range("N81").Value = -506351.5
$ws.Range("H84").Value = 401691.8
$ws.Range("J84").Value = 252114.75
$ws.Range("L84").Value = 2521147.5
$ws.Range("N84").Value = -2531755.5
$ws.Range("H107").Value = 185457
$ws.Range("J107").Value = 337733.34
$ws.Range("L107").Value = 1013200.02
$ws.Range("N107").Value = -1017040.02

Write-Output "Applied updates to ALC/ARM/BSM/CRP/CUL/LTW/WVR sheets"
